$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 767
$ws.Range("F6").Value = 2477
$ws.Range("F8").Value = 1833
$ws.Range("F9").Value = 3143
$ws.Range("F10").Value = 191
$ws.Range("F11").Value = 4632
$ws.Range("F12").Value = 427
$ws.Range("F13").Value = 247
$ws.Range("F15").Value = 589
$ws.Range("F16").Value = 279
$ws.Range("F19").Value = 630
$ws.Range("F20").Value = 267
$ws.Range("F21").Value = 10
$ws.Range("F22").Value = 85
$ws.Range("F23").Value = 129
$ws.Range("F24").Value = 323
$ws.Range("F25").Value = 4619
$ws.Range("F26").Value = 8
$ws.Range("F29").Value = 4988
$ws.Range("F31").Value = 1162
$ws.Range("F33").Value = 633
$ws.Range("F36").Value = 58
$ws.Range("F37").Value = 106
$ws.Range("F38").Value = 742
$ws.Range("F39").Value = 45
$ws.Range("F40").Value = 678
$ws.Range("F41").Value = 672

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 214
$ws.Range("F3").Value = 1061
$ws.Range("F4").Value = 27

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 214
$ws.Range("F4").Value = 1061
$ws.Range("F5").Value = 27
$ws.Range("F8").Value = 767
$ws.Range("F9").Value = 2477
$ws.Range("F11").Value = 1833
$ws.Range("F13").Value = 3143
$ws.Range("F14").Value = 191
$ws.Range("F15").Value = 4632
$ws.Range("F16").Value = 427
$ws.Range("F17").Value = 247
$ws.Range("F19").Value = 589
$ws.Range("F20").Value = 279
$ws.Range("F23").Value = 630
$ws.Range("F24").Value = 267
$ws.Range("F25").Value = 10
$ws.Range("F27").Value = 85
$ws.Range("F28").Value = 129
$ws.Range("F29").Value = 323
$ws.Range("F30").Value = 4619
$ws.Range("F31").Value = 8
$ws.Range("F34").Value = 4988
$ws.Range("F36").Value = 1162
$ws.Range("F38").Value = 633
$ws.Range("F42").Value = 58
$ws.Range("F43").Value = 106
$ws.Range("F44").Value = 742
$ws.Range("F45").Value = 45
$ws.Range("F46").Value = 678
$ws.Range("F47").Value = 672
